$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.126.22"
$ws.Range("E2").Value = "  +1.99%  "
$ws.Range("D3").Value = "1.825.07"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'232.83"
$ws.Range("E5").Value = "  +3.72%  "
$ws.Range("D6").Value = "'0.612"
$ws.Range("E6").Value = "  +1.81%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "'42.21"
$ws.Range("E8").Value = "  +3.27%  "
$ws.Range("E9").Value = "  +6.99%  "
$ws.Range("D10").Value = "'0.0687"
$ws.Range("E10").Value = "  +2.86%  "
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("D12").Value = "2.093.01"
$ws.Range("E12").Value = "  +1.61%  "
$ws.Range("D13").Value = "1.818.58"
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("D14").Value = "'11.14"
$ws.Range("E14").Value = "  +2.95%  "
$ws.Range("D15").Value = "'0.665"
$ws.Range("E15").Value = "  +5.59%  "
$ws.Range("D16").Value = "'4.68"
$ws.Range("D17").Value = "35.103.42"
$ws.Range("E17").Value = "  +1.88%  "
$ws.Range("D18").Value = "'69.91"
$ws.Range("E18").Value = "  +3.79%  "
$ws.Range("D19").Value = "0.0₃0791"
$ws.Range("E19").Value = "  +3.32%  "
$ws.Range("D20").Value = "'239.71"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "'11.80"
$ws.Range("E21").Value = "  +6.86%  "
$ws.Range("E22").Value = "  +12.90%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  +3.48%  "
$ws.Range("D25").Value = "'171.71"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").Value = "'7.76"
$ws.Range("E26").Value = "  +1.39%  "
$ws.Range("D27").Value = "'17.51"
$ws.Range("E27").Value = "  +1.19%  "
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("E29").Value = "  +31.00%  "
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").Value = "3.351.99"
$ws.Range("E31").Value = "  +37.96%  "
$ws.Range("D32").Value = "'0.0554"
$ws.Range("E32").Value = "  +7.99%  "
$ws.Range("D33").Value = "'3.90"
$ws.Range("E33").Value = "  +3.24%  "
$ws.Range("D34").Value = "'3.99"
$ws.Range("E34").Value = "  +4.15%  "
$ws.Range("D35").Value = "'1.78"
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("D36").Value = "'93.20"
$ws.Range("E36").Value = "  +10.93%  "
$ws.Range("D37").Value = "'0.679"
$ws.Range("E37").Value = "  +5.58%  "
$ws.Range("E38").Value = "  +5.78%  "
$ws.Range("D39").Value = "1.322.19"
$ws.Range("E39").Value = "  +0.72%  "
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "'1.28"
$ws.Range("E41").Value = "  +2.27%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'0.991"
$ws.Range("E42").Value = "  +5.93%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'2.34"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'14.73"
$ws.Range("E44").Value = "  -1.56%  "
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("E47").Value = "  +8.00%  "
$ws.Range("D48").Value = "'0.0510"
$ws.Range("E48").Value = "  -1.68%  "
$ws.Range("D49").Value = "2.002.18"
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.0646"
$ws.Range("E51").Value = "  +5.90%  "
